# Insert one new row at row 646, shifting the existing rows 646-687 down to
# 647-688, then populate the newly created row 646 with the new data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("646:646").Insert()

# The date text must stay a plain text string (as all the other dates in
# column A are stored), not get auto-converted into a date serial number.
# Leading the literal with an apostrophe forces text entry like a user
# typing it in Excel; resetting the style back to Normal afterwards drops
# the transient "quote prefix" style so the cell matches its neighbours.
$ws.Range("A646").Value = "'2026/01/15"
$ws.Range("A646").Style = "Normal"

$ws.Range("B646").Value = "木"
$ws.Range("C646").Value = 23
$ws.Range("D646").Value = 201
